# Fixed errors with MSM implementation.
#
# These two sheets store their "number of lines" metric as text-typed cells
# (e.g. "1", "2") even though the text looks numeric. To keep Excel from
# silently re-typing our replacement text as a Number when we assign it, we
# switch the target cells to the Text number format before writing the new
# values - then the round-tripped xlsx keeps them as shared-string ("s")
# cells, matching the original authoring.

$wb = $excel.ActiveWorkbook

# --- Sheet "classNumberOfLines": fix the per-class line counts for three
#     anonymous-inner-class rows (B10:B12). ---------------------------------
$wsClass = $wb.Worksheets.Item("classNumberOfLines")

# NumberFormat is applied cell-by-cell (rather than to the whole B10:B12
# block at once) so every targeted cell reliably keeps the Text format.
foreach ($r in 10..12) {
    $wsClass.Range("B$r").NumberFormat = "@"
}

$wsClass.Range("B10").Value = "0"
$wsClass.Range("B11").Value = "1"
$wsClass.Range("B12").Value = "1"

# --- Sheet "methodNumberOfLines": fix the per-method line counts - every
#     method that was incorrectly reported as "1" line becomes "0". ---------
$wsMethod = $wb.Worksheets.Item("methodNumberOfLines")

$methodRows = @(2,3,5,7,8,9,11,12,14,15,16,19,20,24,25,26,31,32,33,34,38,39,40,41,42,44)

# Same reasoning: format each target cell individually before writing it.
foreach ($r in $methodRows) {
    $wsMethod.Range("C$r").NumberFormat = "@"
}

foreach ($r in $methodRows) {
    $wsMethod.Range("C$r").Value = "0"
}
